$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 and IF in columns I and J, row 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style (bold, bordered, centered) used by the other header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Per-row I0 (inning started) and IF (inning finished) values
$values = @(
    @{Row=2; I0=7; IF=7},
    @{Row=3; I0=7; IF=7},
    @{Row=4; I0=8; IF=8},
    @{Row=5; I0=7; IF=7},
    @{Row=6; I0=8; IF=8},
    @{Row=7; I0=8; IF=8},
    @{Row=8; I0=6; IF=6},
    @{Row=9; I0=7; IF=7},
    @{Row=10; I0=8; IF=9},
    @{Row=11; I0=7; IF=7},
    @{Row=12; I0=5; IF=6},
    @{Row=13; I0=9; IF=9},
    @{Row=14; I0=6; IF=7},
    @{Row=15; I0=7; IF=7},
    @{Row=16; I0=8; IF=8},
    @{Row=17; I0=7; IF=7},
    @{Row=18; I0=8; IF=8},
    @{Row=19; I0=7; IF=7},
    @{Row=20; I0=6; IF=6},
    @{Row=21; I0=7; IF=7},
    @{Row=22; I0=8; IF=8},
    @{Row=23; I0=7; IF=7},
    @{Row=24; I0=7; IF=7},
    @{Row=25; I0=8; IF=8},
    @{Row=26; I0=8; IF=8},
    @{Row=27; I0=7; IF=7},
    @{Row=28; I0=6; IF=7},
    @{Row=29; I0=8; IF=8},
    @{Row=30; I0=8; IF=8},
    @{Row=31; I0=7; IF=7},
    @{Row=32; I0=8; IF=8},
    @{Row=33; I0=7; IF=7},
    @{Row=34; I0=7; IF=7},
    @{Row=35; I0=8; IF=8},
    @{Row=36; I0=7; IF=7},
    @{Row=37; I0=10; IF=10},
    @{Row=38; I0=7; IF=7},
    @{Row=39; I0=7; IF=7},
    @{Row=40; I0=7; IF=7},
    @{Row=41; I0=7; IF=8},
    @{Row=42; I0=10; IF=10},
    @{Row=43; I0=6; IF=7},
    @{Row=44; I0=6; IF=7},
    @{Row=45; I0=6; IF=6},
    @{Row=46; I0=7; IF=7},
    @{Row=47; I0=8; IF=8},
    @{Row=48; I0=8; IF=8},
    @{Row=49; I0=8; IF=8},
    @{Row=50; I0=7; IF=7},
    @{Row=51; I0=8; IF=8},
    @{Row=52; I0=9; IF=9},
    @{Row=53; I0=8; IF=8},
    @{Row=54; I0=7; IF=7},
    @{Row=55; I0=8; IF=8},
    @{Row=56; I0=8; IF=8},
    @{Row=57; I0=7; IF=7},
    @{Row=58; I0=7; IF=7},
    @{Row=59; I0=7; IF=7},
    @{Row=60; I0=8; IF=8},
    @{Row=61; I0=8; IF=8},
    @{Row=62; I0=7; IF=7},
    @{Row=63; I0=5; IF=6},
    @{Row=64; I0=7; IF=7},
    @{Row=65; I0=6; IF=7},
    @{Row=66; I0=11; IF=11},
    @{Row=67; I0=8; IF=8},
    @{Row=68; I0=8; IF=8},
    @{Row=69; I0=9; IF=9},
    @{Row=70; I0=7; IF=7},
    @{Row=71; I0=8; IF=8},
    @{Row=72; I0=6; IF=6},
    @{Row=73; I0=2; IF=2},
    @{Row=74; I0=7; IF=8},
    @{Row=75; I0=9; IF=9},
    @{Row=76; I0=8; IF=8},
    @{Row=77; I0=5; IF=5}
)

foreach ($item in $values) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I0
    $ws.Cells.Item($item.Row, 10).Value = $item.IF
}
